$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value  = 16.18380000000001
$ws.Range("C9").Value  = -10.33900000000001
$ws.Range("E12").Value = 18.32070000000003
$ws.Range("E14").Value = 16.8095
$ws.Range("C18").Value = -12.7927
$ws.Range("C20").Value = -11.50220000000001
$ws.Range("E26").Value = 16.1586
$ws.Range("C27").Value = -12.15989999999999
$ws.Range("E27").Value = 16.72989999999999
$ws.Range("E29").Value = 17.01250000000001
$ws.Range("C35").Value = -11.9953
$ws.Range("E37").Value = 16.58930000000002
$ws.Range("E38").Value = 16.35470000000001
$ws.Range("E51").Value = 17.2581
$ws.Range("E52").Value = 16.935
$ws.Range("E55").Value = 16.3601
$ws.Range("C69").Value = -11.4613
$ws.Range("E69").Value = 16.98910000000002
$ws.Range("E70").Value = 18.05550000000003
$ws.Range("C76").Value = -12.14060000000001
$ws.Range("C78").Value = -11.8507
$ws.Range("E81").Value = 16.7289
$ws.Range("C82").Value = -11.93269999999999
$ws.Range("C83").Value = -14.22999999999999
$ws.Range("E83").Value = 16.2746
$ws.Range("C93").Value = -10.8486
$ws.Range("E102").Value = 16.7732
